# Rebuild the Pandoc/highlighting character styles (SourceCode + *Tok)
# with a light "pygments/tango"-style palette, dropping the dark
# "zenburn"-style 303030 shading that was on every token style.
$d = $word.ActiveDocument

# --- paragraph style: SourceCode -------------------------------------
# Drop the dark <w:shd> fill from the code-block paragraph style while
# keeping its wordWrap=off setting, basedOn and link intact.
$d.Styles("SourceCode").Delete()
$src = $d.Styles.Add("SourceCode", 1)
$src.NameLocal = "Source Code"
$src.BaseStyle = "Normal"
$src.LinkStyle = "VerbatimChar"
$src.ParagraphFormat.WordWrap = $false

# --- character (token) styles ------------------------------------------
$d.Styles("KeywordTok").Delete()
$s = $d.Styles.Add("KeywordTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2125824  # 007020
$s.Font.Bold = $true

$d.Styles("DataTypeTok").Delete()
$s = $d.Styles.Add("DataTypeTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8336  # 902000

$d.Styles("DecValTok").Delete()
$s = $d.Styles.Add("DecValTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 7381056  # 40a070

$d.Styles("BaseNTok").Delete()
$s = $d.Styles.Add("BaseNTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 7381056  # 40a070

$d.Styles("FloatTok").Delete()
$s = $d.Styles.Add("FloatTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 7381056  # 40a070

$d.Styles("ConstantTok").Delete()
$s = $d.Styles.Add("ConstantTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 136  # 880000

$d.Styles("CharTok").Delete()
$s = $d.Styles.Add("CharTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496  # 4070a0

$d.Styles("SpecialCharTok").Delete()
$s = $d.Styles.Add("SpecialCharTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496  # 4070a0

$d.Styles("StringTok").Delete()
$s = $d.Styles.Add("StringTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496  # 4070a0

$d.Styles("VerbatimStringTok").Delete()
$s = $d.Styles.Add("VerbatimStringTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 10514496  # 4070a0

$d.Styles("SpecialStringTok").Delete()
$s = $d.Styles.Add("SpecialStringTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8939195  # bb6688

$d.Styles("ImportTok").Delete()
$s = $d.Styles.Add("ImportTok", 2)
$s.BaseStyle = "VerbatimChar"

$d.Styles("CommentTok").Delete()
$s = $d.Styles.Add("CommentTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392  # 60a0b0
$s.Font.Italic = $true

$d.Styles("DocumentationTok").Delete()
$s = $d.Styles.Add("DocumentationTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2171322  # ba2121
$s.Font.Italic = $true

$d.Styles("AnnotationTok").Delete()
$s = $d.Styles.Add("AnnotationTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392  # 60a0b0
$s.Font.Bold = $true
$s.Font.Italic = $true

$d.Styles("CommentVarTok").Delete()
$s = $d.Styles.Add("CommentVarTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392  # 60a0b0
$s.Font.Bold = $true
$s.Font.Italic = $true

$d.Styles("OtherTok").Delete()
$s = $d.Styles.Add("OtherTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2125824  # 007020

$d.Styles("FunctionTok").Delete()
$s = $d.Styles.Add("FunctionTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8267782  # 06287e

$d.Styles("VariableTok").Delete()
$s = $d.Styles.Add("VariableTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 8132377  # 19177c

$d.Styles("ControlFlowTok").Delete()
$s = $d.Styles.Add("ControlFlowTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2125824  # 007020
$s.Font.Bold = $true

$d.Styles("OperatorTok").Delete()
$s = $d.Styles.Add("OperatorTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 6710886  # 666666

$d.Styles("BuiltInTok").Delete()
$s = $d.Styles.Add("BuiltInTok", 2)
$s.BaseStyle = "VerbatimChar"

$d.Styles("ExtensionTok").Delete()
$s = $d.Styles.Add("ExtensionTok", 2)
$s.BaseStyle = "VerbatimChar"

$d.Styles("PreprocessorTok").Delete()
$s = $d.Styles.Add("PreprocessorTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 31420  # bc7a00

$d.Styles("AttributeTok").Delete()
$s = $d.Styles.Add("AttributeTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 2723965  # 7d9029

$d.Styles("RegionMarkerTok").Delete()
$s = $d.Styles.Add("RegionMarkerTok", 2)
$s.BaseStyle = "VerbatimChar"

$d.Styles("InformationTok").Delete()
$s = $d.Styles.Add("InformationTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392  # 60a0b0
$s.Font.Bold = $true
$s.Font.Italic = $true

$d.Styles("WarningTok").Delete()
$s = $d.Styles.Add("WarningTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 11575392  # 60a0b0
$s.Font.Bold = $true
$s.Font.Italic = $true

$d.Styles("AlertTok").Delete()
$s = $d.Styles.Add("AlertTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 255  # ff0000
$s.Font.Bold = $true

$d.Styles("ErrorTok").Delete()
$s = $d.Styles.Add("ErrorTok", 2)
$s.BaseStyle = "VerbatimChar"
$s.Font.Color = 255  # ff0000
$s.Font.Bold = $true

$d.Styles("NormalTok").Delete()
$s = $d.Styles.Add("NormalTok", 2)
$s.BaseStyle = "VerbatimChar"

Write-Host "Rebuilt SourceCode + 31 token styles (removed dark shading, applied new palette)"
